$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like a plain number (e.g. "240.52").
# Excel.Range.Value auto-converts such strings to numeric doubles, which would
# silently drop meaningful trailing zeros (e.g. "45.20" -> 45.2, "1.00" -> 1).
# Force text storage via the Text number format, write the value, then restore
# the Normal style so the cell keeps its original (unstyled) appearance.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.402.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.985.78'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.04%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.633'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.16%  '
$ws.Range('B7').Value = 'Solana'
$ws.Range('C7').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '56.58'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.86%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '59.28'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.357'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0725'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.12%  '
$ws.Range('E12').Value = '  -6.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.896'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.26'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.278.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.982.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.15'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '35.357.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0833'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '231.28'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.00%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  -5.79%  '
$ws.Range('E26').Value = '  +4.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.10'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.75'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.68%  '
$ws.Range('E30').Value = '  -3.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0583'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.78%  '
$ws.Range('E34').Value = '  +9.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.25'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -10.41%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.25'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.27%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.80'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.85'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('E40').Value = '  -7.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.95%  '
$ws.Range('E42').Value = '  -5.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.08'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.84%  '
$ws.Range('E44').Value = '  -8.42%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '89.88'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.70%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.363.85'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.87'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.20'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.95%  '
